$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine "Nombre" (first name) + "Apellidos" (surname) into a single full name
# in column A, before the surname column gets repurposed.
$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("A3").Value = "Luis López Fernando"
$ws.Range("A4").Value = "Ana Torres Pardo"

# Column B becomes "Localización" with lat/long-style coordinate strings
# (replacing "Apellidos").
$ws.Range("B1").Value = "Localización"
$ws.Range("B2").Value = "32.21 45.34"
$ws.Range("B3").Value = "76.14 98.54"
$ws.Range("B4").Value = "76.34 57.73"

# Column D becomes "Identificador" holding the DNI values that used to live
# in column G (replacing "Fecha nacimiento"). Clear the old date number
# format first so the cells don't keep a stale "m/d/yy"-style format.
$ws.Range("D1").Value = "Identificador"
$ws.Range("D2:D4").ClearFormats()
$ws.Range("D2").Value = "90500084Y"
$ws.Range("D3").Value = "19160962F"
$ws.Range("D4").Value = "09940449X"

# Column E becomes "Tipo" with a constant numeric marker (replacing
# "Dirección postal").
$ws.Range("E1").Value = "Tipo"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1

# Columns F ("Nacionalidad") and G ("DNI") are no longer needed now that the
# identifier moved to column D; drop them entirely.
$ws.Columns("F:G").Delete()

# Move the active selection to where the user left off editing.
$ws.Range("B5").Select()
